$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells in row 1 from "_old" / "_new" suffixes
# to "_FV2310" / "_FV2404" respectively.
$headers = @(
  "Segmentname_FV2310",
  "Segmentgruppe_FV2310",
  "Segment_FV2310",
  "Datenelement_FV2310",
  "Segment ID_FV2310",
  "Code_FV2310",
  "Qualifier_FV2310",
  "Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310",
  "Bedingung_FV2310",
  "diff",
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the data range into a real Excel table (Table1) with an
# autofilter, matching the layout already present in the sheet.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U94"), 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (pane split after row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Done"
